$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 229, pushing existing rows 229:243 down to 230:244
$ws.Rows(229).Insert()

# Populate the new row 229 with the new data record (mostly identical static
# fields to the row that used to be at 229, now at 230, but new date/values)
$ws.Cells.Item(229, 1).Value = 3
$ws.Cells.Item(229, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(229, 3).Value = "Coquimbo"
$ws.Cells.Item(229, 4).Value = 44516
$ws.Cells.Item(229, 4).NumberFormat = $ws.Cells.Item(230, 4).NumberFormat
$ws.Cells.Item(229, 5).Value = 5
$ws.Cells.Item(229, 6).Value = 100112040
$ws.Cells.Item(229, 7).Value = "Cilantro"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 250
$ws.Cells.Item(229, 11).Value = 2000
$ws.Cells.Item(229, 12).Value = 2300
$ws.Cells.Item(229, 13).Value = 2156
$ws.Cells.Item(229, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(229, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(229, 16).Value = 719
$ws.Cells.Item(229, 17).Value = 3
$ws.Cells.Item(229, 18).Value = "Hortaliza"
